# Updated symbol list on Tue Dec 20 07:58:09 UTC 2022 with GitHub Actions
# Refresh the "Price" (column D) and a couple of "Volume(1h)" (column E)
# values to the latest scraped figures. All of these cells hold their
# numbers as literal text (e.g. "0.001590" keeps a trailing zero that a
# real numeric value would drop), so force Text formatting before writing
# so the updated strings are preserved exactly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($cellRef, $value)
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $value
}

Set-TextValue "D2"  "248.37"
Set-TextValue "D3"  "21.68"
Set-TextValue "D4"  "5.296"
Set-TextValue "D6"  "3.428"
Set-TextValue "D7"  "6.377"
Set-TextValue "D8"  "0.8108"
Set-TextValue "D9"  "0.9505"
Set-TextValue "D10" "0.1431"
Set-TextValue "D11" "0.07595"
Set-TextValue "D13" "0.03095"
Set-TextValue "D14" "0.09312"
Set-TextValue "D15" "3.577"
Set-TextValue "D16" "0.001590"
Set-TextValue "D17" "0.04711"
Set-TextValue "D18" "0.0005781"
Set-TextValue "E18" "17OneONEWorstin24h"
Set-TextValue "D19" "0.006298"
Set-TextValue "D20" "0.005055"
Set-TextValue "D22" "0.0001500"
Set-TextValue "D23" "3.784"
Set-TextValue "D24" "2.140"
Set-TextValue "D25" "0.3301"
Set-TextValue "D28" "0.0003000"
Set-TextValue "D40" "0.03957"
Set-TextValue "D41" "0.006840"
Set-TextValue "D42" "0.1065"
Set-TextValue "D43" "0.003030"
Set-TextValue "D45" "0.00005600"
Set-TextValue "D47" "0.0005501"
Set-TextValue "E47" "46ACDXExchangeACXT"
Set-TextValue "D48" "0.7801"
Set-TextValue "D49" "0.1762"
Set-TextValue "D50" "0.00002100"
Set-TextValue "D51" "0.01010"
